# Apply the targets-sheet update: fix the "required_documents" value for the
# existing TEST row and append a new row for "Breega".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: K2 (required_documents) gains the "pitch_deck" option alongside
# the existing "video" one. Typing a brand-new value into the cell (rather
# than just overwriting the shared string) also resets its formatting to the
# column's own default style, so pick up that format from a still-empty cell
# further down column K before writing the new value.
$ws.Range("K500").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("K2").Value = "pitch_deck, video"

# --- Row 3: brand-new "Breega" entry.
$ws.Range("A3").Value = "Breega"
$ws.Range("B3").Value = "https://www.breega.com/"
$ws.Range("C3").Value = "https://rm531z4dws8.typeform.com/to/NNZmuM7H?typeform-source=www.breega.com"
$ws.Range("E3").Value = "form"
$ws.Range("G3").Value = "B2B SaaS, Fintech, Consumer"
$ws.Range("H3").Value = "EMEA, Emerging Markets"
$ws.Range("I3").Value = "standard"
$ws.Range("J3").Value = "11-20"
$ws.Range("K3").Value = "pitch_deck, video"
$ws.Range("L3").Value = "TEST"
$ws.Range("M3").Value = "we back exceptional founders—sometimes before they even see it themselves—from pre-Seed to Series A+, building Digital, Climate, and Deep Tech startups ..."
$ws.Range("N3").Value = "FREE"
$ws.Range("F3").Value = "Pre-seed; Pre-seed; Series A; Series B; Growth"

# --- The K-column dropdown validation used to start at K3 (leaving K2, the
# header data row, unvalidated); now that K2 holds a freshly-typed value it
# is folded back into the validated range.
$found = $false
$dvs = $ws.Range("K2").Validation
try {
  if ($dvs.Type -ne 0) { $found = $true }
} catch {}

$ws.Range("K2:K1000").Validation.Delete()
$ws.Range("K2:K1000").Validation.Add(3, 1, 1, "=ValidationData!$K$1:$K$5")
$ws.Range("K2:K1000").Validation.ErrorTitle = "Guidance"
$ws.Range("K2:K1000").Validation.ErrorMessage = "Select a value from the list, or enter multiple values separated by commas or semicolons."
